$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 76 - new resale number snapshot for 2025-02-18 22:52:00
# Text-like columns (Date/Time/Weekday/Week) are entered with a leading
# apostrophe so they stay literal text (matching the existing inlineStr
# cells) instead of being auto-converted to Excel date/time serials or
# numbers, then the style is reset to "Normal" so no extra number-format
# style gets attached to the cell.
$ws.Range("A76").Value = "'2025-02-18"
$ws.Range("A76").Style = "Normal"

$ws.Range("B76").Value = "'22:52:00"
$ws.Range("B76").Style = "Normal"

$ws.Range("C76").Value = "'Tuesday"
$ws.Range("C76").Style = "Normal"

$ws.Range("D76").Value = "'07"
$ws.Range("D76").Style = "Normal"

$ws.Range("E76").Value = 129089
$ws.Range("F76").Value = 140417
$ws.Range("G76").Value = 171118
$ws.Range("H76").Value = 159545
$ws.Range("I76").Value = -1
$ws.Range("J76").Value = 145728
$ws.Range("K76").Value = -1
$ws.Range("L76").Value = -1
$ws.Range("M76").Value = 192220
$ws.Range("N76").Value = 114953
$ws.Range("O76").Value = 45595
$ws.Range("P76").Value = 28996
$ws.Range("Q76").Value = 67075
$ws.Range("R76").Value = -1
$ws.Range("S76").Value = 47100
$ws.Range("T76").Value = -1
